# Revert "Epochs = 700" back to "Epochs = 600" in the clusternetasmodel notes.
#
# The target OOXML keeps the paragraph's own run (with the original
# "Epochs = " + "00 " text) but re-types the single changed digit ("6")
# as its own run, exactly like Word does when a user selects one
# character and types its replacement - the run gets split into
# "before" / "typed" / "after" pieces instead of being merged back into
# one <w:r>. We reproduce that run layout explicitly.

$d = $word.ActiveDocument
$target = $d.Content

$found = $target.Find.Execute("Epochs = 700 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Epochs = 700 ' in the document"
}

# Capture the enclosing paragraph's real OOXML (pPr, rsids, paraId, ...)
# so we can splice in the new runs without inventing/guessing any
# paragraph-level attributes.
$full = $target.WordOpenXML

$needle = "Epochs = 700"
$needleIdx = $full.IndexOf($needle)

$prefix = $full.Substring(0, $needleIdx)
$pOpenIdx = $prefix.LastIndexOf("<w:p ")
$pOpenIdxAlt = $prefix.LastIndexOf("<w:p>")
if ($pOpenIdxAlt -gt $pOpenIdx) { $pOpenIdx = $pOpenIdxAlt }

$pCloseIdx = $full.IndexOf("</w:p>", $needleIdx)
$paraXml = $full.Substring($pOpenIdx, $pCloseIdx - $pOpenIdx + 6)

$oldRun = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Epochs = 700 </w:t></w:r>'
$newRuns = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Epochs = </w:t></w:r>' `
    + '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>6</w:t></w:r>' `
    + '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">00 </w:t></w:r>'

$newParaXml = $paraXml -replace [regex]::Escape($oldRun), $newRuns
if ($newParaXml -eq $paraXml) {
    throw "Did not find the expected single run inside the target paragraph"
}

# Replace just that paragraph's contents: delete the old run's text, then
# insert the rebuilt paragraph XML (with the split runs) at that spot.
$target.Delete()
$target.Collapse(1)

$package = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' `
    + '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($package)
